$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.944.22'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '3.335.46'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '584.91'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = '177.88'
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("E9").Value = '  +4.27%  '
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("D11").Value = '48.10'
$ws.Range("E11").Value = '  +6.06%  '
$ws.Range("E12").Value = '  +1.78%  '
$ws.Range("D13").Value = '698.08'
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("D14").Value = '3.875.57'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '8.45'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").Value = '67.990.54'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D18").Value = '3.350.83'
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("D20").Value = '11.17'
$ws.Range("E20").Value = '  +2.80%  '
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").Value = '5.44'
$ws.Range("E22").Value = '  +1.62%  '
$ws.Range("D23").Value = '16.92'
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("D24").Value = '100.23'
$ws.Range("E24").Value = '  +2.93%  '
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = '9.48'
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("D28").Value = '33.06'
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("D29").Value = '8.57'
$ws.Range("E29").Value = '  +2.00%  '
$ws.Range("E30").Value = '  -4.38%  '
$ws.Range("D31").Value = '579.22'
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("E33").Value = '  +1.82%  '
$ws.Range("D34").Value = '3.740.61'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").Value = '57.39'
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '3.36'
$ws.Range("E37").Value = '  +1.36%  '
$ws.Range("D38").Value = '35.42'
$ws.Range("E38").Value = '  +9.02%  '
$ws.Range("D39").Value = '0.135'
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = '2.63'
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("D42").Value = '0.0₃0675'
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E49").Value = '  -1.31%  '
$ws.Range("D50").Value = '130.86'
$ws.Range("E50").Value = '  +2.62%  '
$ws.Range("D51").Value = '2.62'
$ws.Range("E51").Value = '  -0.06%  '
